$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 71

$ws.Cells.Item($row, 1).Value = 46020
$ws.Cells.Item($row, 2).Value = 158
$ws.Cells.Item($row, 3).Value = 167
$ws.Cells.Item($row, 4).Value = 158

# Match the formatting of the cell above (date-style column A)
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
